$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case WAT43 goes into row 30. The "JIRA ID" column (B) value for
# the existing row 29 (WAT42) is also corrected from WAT-194 to WAT-207,
# and the new row 30 gets WAT-197 in that column.
# Order of writes matters for the shared-strings table layout, so keep it
# matching the order the strings were authored in.
$ws.Range("A30").Value = "WAT43"
$ws.Range("B29").Value = "WAT-207"
$ws.Range("C30").Value = "Verify that, If system retrieves only one country and only one organization`nSystem must directly display the results in search results page if the search result count is less than 50."
$ws.Range("B30").Value = "WAT-197"
$ws.Range("D30").Value = "Y"

# Row 29 is the template for formatting (borders / wrap-text / etc.) of the
# new row, same as every other data row in this sheet.
$ws.Range("A29:E29").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(30).RowHeight = 30

$ws.Range("B30").Select()
